$d = $word.ActiveDocument

$d.Content.Find.Execute("The line numbers we refer to below are taken from the “clean” version of the manuscript we have submitted (manuscript.pdf).", $true, $false, $false, $false, $false, $true, 1, $false, "The line numbers we refer to below are taken from the “clean” version of the manuscript we have submitted (manuscript.docx).", 2) | Out-Null
$d.Content.Find.Execute("Thank you for highlighting this area that would benefit from further clarification. We decided to aggregate the data from each antibiotic across different levels because we were hoping to capture the common differences associated with C. difficile colonization as opposed to the effects of the antibiotic. To be more transparent in this decision, we added a sentence to explain the aggregation (Lines 106-109). Also to ensure that the differences we observed were not merely differences due to antibiotic dosages, we added a supplemental figure (Figure S3) to show the alpha diversity by dose of cefoperazone. None of the comparisons were statistically significant, which we added a sentence to state this (Lines 147-148).", $true, $false, $false, $false, $false, $true, 1, $false, "Thank you for highlighting this area that would benefit from further clarification. We decided to aggregate the data from each antibiotic across different levels because we were hoping to capture the common differences associated with C. difficile colonization as opposed to the effects of the antibiotic. To be more transparent in this decision, we added a sentence to explain the aggregation (Lines 106-109). Also to ensure that the differences we observed were not merely differences due to antibiotic dosages, we added a supplemental figure (Figure S3) to show the alpha diversity by dose of cefoperazone. None of the comparisons were statistically significant, which we added a sentence to state this (Lines 146-147).", 2) | Out-Null
$d.Content.Find.Execute("We acknowledge that individual models would have allowed us to identify the OTUs most predictive of each treatment. While we treated them separately in previous comparisons, here with machine learning we felt it was best to model the data together. This way it would provide sufficient sample size to create a robust model and not overfit to one specific treatment. We have added a sentence to state our reasoning for modeling all treatments together (Lines 202-204)", $true, $false, $false, $false, $false, $true, 1, $false, "We acknowledge that individual models would have allowed us to identify the OTUs most predictive of each treatment. While we treated them separately in previous comparisons, here with machine learning we felt it was best to model the data together. This way it would provide sufficient sample size to create a robust model and not overfit to one specific treatment. We have added a sentence to state our reasoning for modeling all treatments together (Lines 201-203)", 2) | Out-Null
$d.Content.Find.Execute("t we were implying by “greater effects” with examples such as the one recommended by this comment by stating the more virulent strains drive a stronger immune response, increasing inflammation and inflammation-associated bacteria such as Enterobacteriaceae (Lines 324-329).", $true, $false, $false, $false, $false, $true, 1, $false, "t we were implying by “greater effects” with examples such as the one recommended by this comment by stating the more virulent strains drive a stronger immune response, increasing inflammation and inflammation-associated bacteria such as Enterobacteriaceae (Lines 324-328).", 2) | Out-Null
$d.Content.Find.Execute("We agree with you and have made this change by describing it as a “subset” (Line 131).", $true, $false, $false, $false, $false, $true, 1, $false, "We agree with you and have made this change by describing it as a “subset” (Line 130).", 2) | Out-Null
$d.Content.Find.Execute("Thank you for your suggestion. We understand that using Figure 1E isn’t effective at demonstrating what this statement was describing. We have incorporated your comment by creating a supplemental plot (Figure S2) to show the distribution of OTU relative abundance in cefoperazone treated mice. The plot shows that the mice with increased alpha diversity have an increase in abundance of otus in lower abundance in all other samples. We revised the statement as well to describe this (Lines 142-145).", $true, $false, $false, $false, $false, $true, 1, $false, "Thank you for your suggestion. We understand that using Figure 1E isn’t effective at demonstrating what this statement was describing. We have incorporated your comment by creating a supplemental plot (Figure S2) to show the distribution of OTU relative abundance in cefoperazone treated mice. The plot shows that the mice with increased alpha diversity have an increase in abundance of otus in lower abundance in all other samples. We revised the statement as well to describe this (Lines 141-143).", 2) | Out-Null
$d.Content.Find.Execute("We have reflected this comment by more explicitly describing the details and aspects of the referenced articles to state that these articles described bacteria associated with C. difficile colonization, and then described the outcomes when trying to use identified bacteria to affect C. difficile colonization (Lines 270-275).", $true, $false, $false, $false, $false, $true, 1, $false, "We have reflected this comment by more explicitly describing the details and aspects of the referenced articles to state that these articles described bacteria associated with C. difficile colonization, and then described the outcomes when trying to use identified bacteria to affect C. difficile colonization (Lines 269-274).", 2) | Out-Null
$d.Content.Find.Execute("This comment points out the number of differences in utilized substrates reported in previous studies in our lab. However, we believe niches available are not limited to the number of unique substrates metabolized, but also includes the amount of resources and the physical space that is available. We added text to help clarify this meaning (Lines 299-301)", $true, $false, $false, $false, $false, $true, 1, $false, "This comment points out the number of differences in utilized substrates reported in previous studies in our lab. However, we believe niches available are not limited to the number of unique substrates metabolized, but also includes the amount of resources and the physical space that is available. We added text to help clarify this meaning (Lines 298-300)", 2) | Out-Null
$d.Content.Find.Execute("Thank you for your perspective, however we are unsure what this comment is requesting to be revised. We state that it is possible that an FMT may not be sufficient to recover a significantly disrupted microbiome, such as our hypothesis for the cefoperzaone-treated mice. It seems like this comment is stating since Seekatz et al 2015 showed FMT prevented relapse with cefoperazone that our hypothesis is incorrect. However, Seekatz et al 2015 have a much different experimental scenario, which used a different breeding colony and sensitized using cefoperzaone, challenged with C. difficile 630, 4 days later given a 5 day course of vancomycin, then two days later given two daily gavages of an FMT. After that there was a transient increase in CFU which remained below LOD until the mice were given an IP of clindamycin which caused another transient increase in CFU. So I feel the referenced manuscript is not directly applicable to lines 273-276. Additionally, the comment states the rarity of FMT failure (Seekatz et al. 2014).  While an FMT is ~80-90% effective, this still leaves 10-20% of patients in which it does not work. So even if it is relatively rare, it still occurs. ", $true, $false, $false, $false, $false, $true, 1, $false, "Thank you for your perspective, however we are unsure what this comment is requesting to be revised. We state that it is possible that an FMT may not be sufficient to recover a significantly disrupted microbiome, such as our hypothesis for the cefoperzaone-treated mice. It seems like this comment is stating since Seekatz et al 2015 showed FMT prevented relapse with cefoperazone that our hypothesis is incorrect. However, Seekatz et al 2015 have a much different experimental scenario, which used a different breeding colony and sensitized using cefoperzaone, challenged with C. difficile 630, 4 days later given a 5 day course of vancomycin, then two days later given two daily gavages of an FMT. After that there was a transient increase in CFU which remained below LOD until the mice were given an IP of clindamycin which caused another transient increase in CFU. So I feel the referenced manuscript is not directly applicable to lines 272-274. Additionally, the comment states the rarity of FMT failure (Seekatz et al. 2014).  While an FMT is ~80-90% effective, this still leaves 10-20% of patients in which it does not work. So even if it is relatively rare, it still occurs. ", 2) | Out-Null
$d.Content.Find.Execute("Thank you for your suggestion. While the Jenior et al. experiments are similar to the ones here, there are differences that make direct application of those findings to these difficult. Those experiments looked at activity at peak infection, whereas we are looking at the taxonomic changes associated with clearance. There may be overlap but there also could be differences taxonomically as well as nutritionally through the clearance process. Also this comment associates phylogenetic diversity with functional diversity. We acknowledge that we were not clear about the limited disruption. We believe it is possible that specific OTUs have specific functions not shared by other species/OTUs of the same genus. We clarify that we believe, based on the changes we observed in figures 3 and 4, that only a few specific bacteria are necessary to clear colonization in those communities (lines 309-311, 335-337).", $true, $false, $false, $false, $false, $true, 1, $false, "Thank you for your suggestion. While the Jenior et al. experiments are similar to the ones here, there are differences that make direct application of those findings to these difficult. Those experiments looked at activity at peak infection, whereas we are looking at the taxonomic changes associated with clearance. There may be overlap but there also could be differences taxonomically as well as nutritionally through the clearance process. Also this comment associates phylogenetic diversity with functional diversity. We acknowledge that we were not clear about the limited disruption. We believe it is possible that specific OTUs have specific functions not shared by other species/OTUs of the same genus. We clarify that we believe, based on the changes we observed in figures 3 and 4, that only a few specific bacteria are necessary to clear colonization in those communities (lines 308-310, 334-336).", 2) | Out-Null
$d.Content.Find.Execute("We agree that description of the points may not be correctly interpreted so we added to their descriptions to ensure they will be correctly interpreted  “Dark larger points in foreground are median relative abundance and transparent smaller points in background are relative abundance of individual mice.” (Lines 604-605). Also we have completed the first sentence with “were identified” (Line 603-604)", $true, $false, $false, $false, $false, $true, 1, $false, "We agree that description of the points may not be correctly interpreted so we added to their descriptions to ensure they will be correctly interpreted  “Dark larger points in foreground are median relative abundance and transparent smaller points in background are relative abundance of individual mice.” (Lines 600-601). Also we have completed the first sentence with “were identified” (Line 599-600)", 2) | Out-Null
$d.Content.Find.Execute("We agree with you, it was not clear how to differentiate or identify the lines mentioned in the legend. We changed IQR to a light green band and AUROC to a dark green line to make them more easily differentiated from each other and boxplots and updated figure legend to match. (Line 627-629)", $true, $false, $false, $false, $false, $true, 1, $false, "We agree with you, it was not clear how to differentiate or identify the lines mentioned in the legend. We changed IQR to a light green band and AUROC to a dark green line to make them more easily differentiated from each other and boxplots and updated figure legend to match. (Line 619-621)", 2) | Out-Null
$d.Content.Find.Execute("We corrected the “colonization” to “colonized” (Line 185).", $true, $false, $false, $false, $false, $true, 1, $false, "We corrected the “colonization” to “colonized” (Line 184).", 2) | Out-Null
$d.Content.Find.Execute("Thank you for identifying these issues in our references. We edited the reference input to ensure proper formatting for all of our references to have appropriate spacing, capitalization and italics. (Lines 431, 434, 474)", $true, $false, $false, $false, $false, $true, 1, $false, "Thank you for identifying these issues in our references. We edited the reference input to ensure proper formatting for all of our references to have appropriate spacing, capitalization and italics. (Lines 430, 435, 474)", 2) | Out-Null
$d.Content.Find.Execute("Lines 335-339", $true, $false, $false, $false, $false, $true, 1, $false, "Lines 335-338", 2) | Out-Null
